# daily auto push: 2026-01-05 09:44 UTC
# Insert one new row of data (2026/01/05, 月, 13, 201) right before the
# existing "2026/12/29" block, shifting all rows from 571 downward by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 571; everything currently at/after
# row 571 (the 2026/12/29 ... 2027/01/05 block) shifts down to 572..613.
$ws.Rows.Item(571).Insert()

$dateCell = $ws.Cells.Item(571, 1)
# Force text storage so "2026/01/05" is kept as a literal string (matching
# every other date cell in the column) instead of being auto-parsed into a
# date serial number by the COM layer.
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/01/05"
# Drop back to the workbook's default style so the cell doesn't end up
# carrying a stray "Text" number-format style like its neighbours.
$dateCell.Style = "Normal"

$ws.Cells.Item(571, 2).Value = "月"
$ws.Cells.Item(571, 3).Value = 13
$ws.Cells.Item(571, 4).Value = 201
